# Termino el switch del lectorExcel
# The reader used to leave the A3:A4 cells merged with a single label in A3
# and a blank A4. The switch-statement rewrite now writes the label into
# every row of the activity instead of relying on the merge, so the merge
# is removed and A4 gets its own copy of the "Combustion fija" text (with
# the same "no special alignment" styling the lector applies to those
# cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The activity column (A3:A4) is no longer merged - each row now carries
# its own value.
$ws.Range("A3:A4").UnMerge()

# A4 gets the same activity label as A3 ("Combustion fija").
$ws.Range("A4").Value = $ws.Range("A3").Value()

# Both cells pick up the (new) plain style the lector applies once it no
# longer needs the centered merged-cell look.
$ws.Range("A3").HorizontalAlignment = 1
$ws.Range("A4").HorizontalAlignment = 1

# Leftover selection state from the editing session that produced the file.
$null = $ws.Range("E11").Select()
